# Update "Training Dashboard" sheet with new progress as of 04-Nov-2025.
# Row 3: PERIOD TO EXPIRE 477 -> 476, LAST UPDATE 03-Nov-2025 -> 04-Nov-2025
# Row 4: PERIOD TO EXPIRE 700 -> 699, LAST UPDATE 03-Nov-2025 -> 04-Nov-2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$ws.Range("H3").Value = 476
$ws.Range("H4").Value = 699

# I3/I4 hold the "LAST UPDATE" date but are stored as plain text, not real
# Excel dates. Assigning a date-looking string straight to .Value would make
# Excel auto-convert it into a date serial number, so instead write it as a
# text formula and then collapse the formula down to its static result via
# copy / paste-special-values - this keeps the cell as literal text (and
# keeps the existing cell style) instead of turning it into a date.
$i3 = $ws.Range("I3")
$i3.Formula = "=""04-Nov-2025"""
$i3.Copy()
$i3.PasteSpecial(-4163)

$i4 = $ws.Range("I4")
$i4.Formula = "=""04-Nov-2025"""
$i4.Copy()
$i4.PasteSpecial(-4163)

$excel.CutCopyMode = $false
